# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row (by column F) -> new value
$updates = @{
    3 = 2515
    4 = 491
    6 = 6542
    7 = 378
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
